$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: logistic_embeddings
$ws.Range("C5").Value = 0.333
$ws.Range("D5").Value = 0.427
$ws.Range("E5").Value = 0.447
$ws.Range("F5").Value = 0.479
$ws.Range("G5").Value = 0.46
$ws.Range("H5").Value = 0.484

# Row 7: classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.333
$ws.Range("D7").Value = 0.427

# Row 8: BERT-base
$ws.Range("C8").Value = 0.328
$ws.Range("D8").Value = 0.531
$ws.Range("E8").Value = 0.554
$ws.Range("F8").Value = 0.607
$ws.Range("G8").Value = 0.591
$ws.Range("H8").Value = 0.618

# Row 9: BERT-base-nli
$ws.Range("B9").Value = 0.359
$ws.Range("C9").Value = 0.521
$ws.Range("D9").Value = 0.615
$ws.Range("E9").Value = 0.631
$ws.Range("F9").Value = 0.673
$ws.Range("G9").Value = 0.63
$ws.Range("H9").Value = 0.648

$wb.Save()
